$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Solver Options" to "Sheet1"
$ws.Name = "Sheet1"

# New column C header, styled like the existing header cells (A1/B1):
# bold font, thin border all around, centered horizontally, top vertically.
$ws.Cells.Item(1, 3).Value = "Count"
$ws.Cells.Item(1, 3).Font.Bold = $true
$ws.Cells.Item(1, 3).HorizontalAlignment = -4108
$ws.Cells.Item(1, 3).VerticalAlignment = -4160
$ws.Cells.Item(1, 3).Borders.LineStyle = 1

# Rows that carry a real (numeric) count of 0; everything else in the
# new column is a blank/empty entry (present in the column, but no value).
$zeroRows = @(13, 19, 29)

for ($r = 2; $r -le 37; $r++) {
    if ($zeroRows -contains $r) {
        $ws.Cells.Item($r, 3).Value = 0
    } else {
        # Force the cell to exist (blank) without deleting it the way an
        # empty-string Value assignment would: copy formatting from an
        # already-present, unstyled cell in the same row so no content is
        # written but the cell is materialized in the sheet.
        $ws.Cells.Item($r, 3).Style = $ws.Cells.Item($r, 1).Style
    }
}
